# Auto-generated edit script: updates FFXIV leve-profit price/profit columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ /
# LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across several sheets, per refreshed
# market-board data from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2998.3333
$ws.Range("I62").Value = 2997
$ws.Range("J62").Value = 2999
$ws.Range("K62").Value = 2997
$ws.Range("L62").Value = 2999
$ws.Range("M62").Value = -2373
$ws.Range("N62").Value = -4247
$ws.Range("H65").Value = 2998.3333
$ws.Range("I65").Value = 2997
$ws.Range("J65").Value = 2999
$ws.Range("K65").Value = 14985
$ws.Range("L65").Value = 14995
$ws.Range("M65").Value = -11865
$ws.Range("N65").Value = -21235
$ws.Range("H74").Value = 2900
$ws.Range("I74").Value = 2900
$ws.Range("K74").Value = 2900
$ws.Range("M74").Value = -1964
$ws.Range("H76").Value = 3910086.2
$ws.Range("I76").Value = 7811641
$ws.Range("K76").Value = 7811641
$ws.Range("M76").Value = -7811326
$ws.Range("H77").Value = 2900
$ws.Range("I77").Value = 2900
$ws.Range("K77").Value = 14500
$ws.Range("M77").Value = -9820
$ws.Range("H79").Value = 3910086.2
$ws.Range("I79").Value = 7811641
$ws.Range("K79").Value = 7811641
$ws.Range("M79").Value = -7810549
$ws.Range("H106").Value = 2449.25
$ws.Range("I106").Value = 2449.25
$ws.Range("K106").Value = 2449.25
$ws.Range("M106").Value = -1818.25
$ws.Range("H112").Value = 5813.4287
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 5813.4287
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 17440.2861
$ws.Range("N112").Value = -19656.2861
$ws.Range("H135").Value = 564.58826
$ws.Range("I135").Value = 555.13336
$ws.Range("J135").Value = 635.5
$ws.Range("K135").Value = 4996.20024
$ws.Range("L135").Value = 5719.5
$ws.Range("M135").Value = -2461.20024
$ws.Range("N135").Value = -10789.5
$ws.Range("H138").Value = 2873.0698
$ws.Range("I138").Value = 2498.9285
$ws.Range("J138").Value = 3571.4666
$ws.Range("K138").Value = 7496.7855
$ws.Range("L138").Value = 10714.3998
$ws.Range("M138").Value = -2356.7855
$ws.Range("N138").Value = -20994.3998
$ws.Range("M112").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 117.5
$ws.Range("I4").Value = 117.5
$ws.Range("K4").Value = 117.5
$ws.Range("M4").Value = -1.5
$ws.Range("H32").Value = 2806.3809
$ws.Range("I32").Value = 2160.6572
$ws.Range("J32").Value = 6035
$ws.Range("K32").Value = 2160.6572
$ws.Range("L32").Value = 6035
$ws.Range("M32").Value = -1873.6572
$ws.Range("N32").Value = -6609
$ws.Range("H38").Value = 519
$ws.Range("I38").Value = 519
$ws.Range("K38").Value = 519
$ws.Range("M38").Value = -52
$ws.Range("H74").Value = 1184.225
$ws.Range("I74").Value = 958.9666999999999
$ws.Range("K74").Value = 958.9666999999999
$ws.Range("M74").Value = -84.96669999999995
$ws.Range("H77").Value = 1184.225
$ws.Range("I77").Value = 958.9666999999999
$ws.Range("K77").Value = 4794.8335
$ws.Range("M77").Value = -426.8334999999997
$ws.Range("H97").Value = 595.93335
$ws.Range("I97").Value = 412.72726
$ws.Range("K97").Value = 412.72726
$ws.Range("M97").Value = 83.27274
$ws.Range("H122").Value = 2828.625
$ws.Range("I122").Value = 1804.1428
$ws.Range("K122").Value = 5412.428400000001
$ws.Range("M122").Value = -2962.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 925
$ws.Range("I22").Value = 850
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 850
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -677
$ws.Range("N22").Value = -1346
$ws.Range("H80").Value = 8693.916999999999
$ws.Range("J80").Value = 11475.556
$ws.Range("L80").Value = 11475.556
$ws.Range("N80").Value = -13471.556
$ws.Range("H83").Value = 8693.916999999999
$ws.Range("J83").Value = 11475.556
$ws.Range("L83").Value = 57377.78
$ws.Range("N83").Value = -67361.78
$ws.Range("H94").Value = 468.05405
$ws.Range("I94").Value = 343.87878
$ws.Range("J94").Value = 1492.5
$ws.Range("K94").Value = 343.87878
$ws.Range("L94").Value = 1492.5
$ws.Range("M94").Value = 107.12122
$ws.Range("N94").Value = -2394.5
$ws.Range("H105").Value = 2161.963
$ws.Range("I105").Value = 2119.75
$ws.Range("J105").Value = 2499.6667
$ws.Range("K105").Value = 2119.75
$ws.Range("L105").Value = 2499.6667
$ws.Range("M105").Value = -372.75
$ws.Range("N105").Value = -5993.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 854.3
$ws.Range("I16").Value = 818.25
$ws.Range("J16").Value = 998.5
$ws.Range("K16").Value = 818.25
$ws.Range("L16").Value = 998.5
$ws.Range("M16").Value = -531.25
$ws.Range("N16").Value = -1572.5
$ws.Range("H31").Value = 1872.5
$ws.Range("I31").Value = 1423.2
$ws.Range("J31").Value = 2621.3333
$ws.Range("K31").Value = 1423.2
$ws.Range("L31").Value = 2621.3333
$ws.Range("M31").Value = -1128.2
$ws.Range("N31").Value = -3211.3333
$ws.Range("H34").Value = 1872.5
$ws.Range("I34").Value = 1423.2
$ws.Range("J34").Value = 2621.3333
$ws.Range("K34").Value = 1423.2
$ws.Range("L34").Value = 2621.3333
$ws.Range("M34").Value = -1221.2
$ws.Range("N34").Value = -3025.3333
$ws.Range("H58").Value = 1612172.2
$ws.Range("I58").Value = 2289606.8
$ws.Range("K58").Value = 2289606.8
$ws.Range("M58").Value = -2289403.8
$ws.Range("H105").Value = 1240
$ws.Range("I105").Value = 1240
$ws.Range("K105").Value = 1240
$ws.Range("M105").Value = 507
$ws.Range("H107").Value = 343.6
$ws.Range("I107").Value = 338.30435
$ws.Range("J107").Value = 353.75
$ws.Range("K107").Value = 338.30435
$ws.Range("L107").Value = 353.75
$ws.Range("M107").Value = 1581.69565
$ws.Range("N107").Value = -4193.75
$ws.Range("H113").Value = 854.3
$ws.Range("I113").Value = 818.25
$ws.Range("J113").Value = 998.5
$ws.Range("K113").Value = 818.25
$ws.Range("L113").Value = 998.5
$ws.Range("M113").Value = 1351.75
$ws.Range("N113").Value = -5338.5
$ws.Range("H132").Value = 2902.5557
$ws.Range("I132").Value = 1483.4
$ws.Range("K132").Value = 4450.200000000001
$ws.Range("M132").Value = -1920.200000000001
$ws.Range("H134").Value = 1120.3889
$ws.Range("I134").Value = 1138.909
$ws.Range("J134").Value = 916.6667
$ws.Range("K134").Value = 3416.727
$ws.Range("L134").Value = 2750.0001
$ws.Range("M134").Value = -881.7270000000003
$ws.Range("N134").Value = -7820.0001
$ws.Range("H136").Value = 1612172.2
$ws.Range("I136").Value = 2289606.8
$ws.Range("K136").Value = 6868820.399999999
$ws.Range("M136").Value = -6866270.399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 658.0417
$ws.Range("J107").Value = 819.64703
$ws.Range("L107").Value = 2458.94109
$ws.Range("N107").Value = -6298.94109
$ws.Range("H122").Value = 1012.1
$ws.Range("I122").Value = 585
$ws.Range("K122").Value = 5265
$ws.Range("M122").Value = -2815
$ws.Range("H131").Value = 809.54
$ws.Range("J131").Value = 816.0205999999999
$ws.Range("L131").Value = 2448.0618
$ws.Range("N131").Value = -12528.0618

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14287.5
$ws.Range("I70").Value = 23950
$ws.Range("K70").Value = 23950
$ws.Range("M70").Value = -23680
$ws.Range("H73").Value = 14287.5
$ws.Range("I73").Value = 23950
$ws.Range("K73").Value = 23950
$ws.Range("M73").Value = -23014
$ws.Range("H122").Value = 1721
$ws.Range("I122").Value = 1335.25
$ws.Range("J122").Value = 2299.625
$ws.Range("K122").Value = 4005.75
$ws.Range("L122").Value = 6898.875
$ws.Range("M122").Value = -1555.75
$ws.Range("N122").Value = -11798.875
$ws.Range("H132").Value = 1328015.1
$ws.Range("I132").Value = 1832857.5
$ws.Range("J132").Value = 2804.125
$ws.Range("K132").Value = 5498572.5
$ws.Range("L132").Value = 8412.375
$ws.Range("M132").Value = -5496042.5
$ws.Range("N132").Value = -13472.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 6495.25
$ws.Range("I32").Value = 5370.3335
$ws.Range("K32").Value = 5370.3335
$ws.Range("M32").Value = -5053.3335
$ws.Range("H55").Value = 334.0645
$ws.Range("J55").Value = 474.8
$ws.Range("L55").Value = 474.8
$ws.Range("N55").Value = -820.8
$ws.Range("H82").Value = 2039.375
$ws.Range("I82").Value = 1469.5
$ws.Range("K82").Value = 1469.5
$ws.Range("M82").Value = -1108.5
$ws.Range("H85").Value = 2039.375
$ws.Range("I85").Value = 1469.5
$ws.Range("K85").Value = 1469.5
$ws.Range("M85").Value = -221.5
$ws.Range("H93").Value = 800
$ws.Range("I93").Value = 700
$ws.Range("K93").Value = 700
$ws.Range("M93").Value = 548
$ws.Range("H100").Value = 1197.8
$ws.Range("J100").Value = 1099.6666
$ws.Range("L100").Value = 1099.6666
$ws.Range("N100").Value = -2181.6666
$ws.Range("H122").Value = 6316.3335
$ws.Range("I122").Value = 1632.6666
$ws.Range("K122").Value = 4897.9998
$ws.Range("M122").Value = -2447.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 61499.92
$ws.Range("I122").Value = 112253.57
$ws.Range("K122").Value = 336760.71
$ws.Range("M122").Value = -334310.71
$ws.Range("H126").Value = 12375.083
$ws.Range("J126").Value = 7875
$ws.Range("L126").Value = 23625
$ws.Range("N126").Value = -28565

Write-Output "Applied Tonberry_Profits market data refresh."
